# Update cryptocurrency price/volume data to latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '26.446.65'
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.31%  '
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.834.94'
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.66%  '
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.21%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '260.27'
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.49%  '
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.22%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.5328'
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.81%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.3009'
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -6.68%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.06861'
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.85%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '17.71'
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -5.78%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.852.17'
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.42%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.7357'
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -5.69%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07299'
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -5.99%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '88.90'
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.39%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.962'
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.24%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.18%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.90'
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.55%  '
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.17%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.000007895'
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.84%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '26.483.59'
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.24%  '
$cell = $ws.Cells.Item(21, 2)
$cell.NumberFormat = "@"
$cell.Value = 'WrappedliquidstakedEther2.0'
$cell = $ws.Cells.Item(21, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.084.29'
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.55%  '
$cell = $ws.Cells.Item(22, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Uniswap'
$cell = $ws.Cells.Item(22, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.579'
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.16%  '
$cell = $ws.Cells.Item(23, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Chainlink'
$cell = $ws.Cells.Item(23, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.962'
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.75%  '
$cell = $ws.Cells.Item(24, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Cosmos'
$cell = $ws.Cells.Item(24, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.245'
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.38%  '
$cell = $ws.Cells.Item(25, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Monero'
$cell = $ws.Cells.Item(25, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '142.64'
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.27%  '
$cell = $ws.Cells.Item(26, 2)
$cell.NumberFormat = "@"
$cell.Value = 'LidoDAOToken'
$cell = $ws.Cells.Item(26, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.212'
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.69%  '
$cell = $ws.Cells.Item(27, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Toncoin'
$cell = $ws.Cells.Item(27, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.683'
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.04%  '
$cell = $ws.Cells.Item(28, 2)
$cell.NumberFormat = "@"
$cell.Value = 'EthereumClassic'
$cell = $ws.Cells.Item(28, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.93'
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.57%  '
$cell = $ws.Cells.Item(29, 2)
$cell.NumberFormat = "@"
$cell.Value = 'BitcoinCash'
$cell = $ws.Cells.Item(29, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '110.23'
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.27%  '
$cell = $ws.Cells.Item(30, 2)
$cell.NumberFormat = "@"
$cell.Value = 'InternetComputer(DFINITY)'
$cell = $ws.Cells.Item(30, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.229'
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.12%  '
$cell = $ws.Cells.Item(31, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Stellar'
$cell = $ws.Cells.Item(31, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.08808'
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.71%  '
$cell = $ws.Cells.Item(32, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Filecoin'
$cell = $ws.Cells.Item(32, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.018'
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.53%  '
$cell = $ws.Cells.Item(33, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Hedera'
$cell = $ws.Cells.Item(33, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.04794'
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.86%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.940'
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.56%  '
$cell = $ws.Cells.Item(35, 2)
$cell.NumberFormat = "@"
$cell.Value = 'ImmutableX'
$cell = $ws.Cells.Item(35, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.7308'
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.52%  '
$cell = $ws.Cells.Item(36, 2)
$cell.NumberFormat = "@"
$cell.Value = 'ARBITRUM'
$cell = $ws.Cells.Item(36, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.129'
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.10%  '
$cell = $ws.Cells.Item(37, 2)
$cell.NumberFormat = "@"
$cell.Value = 'MXToken'
$cell = $ws.Cells.Item(37, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.092'
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.50%  '
$cell = $ws.Cells.Item(38, 2)
$cell.NumberFormat = "@"
$cell.Value = 'RenderToken'
$cell = $ws.Cells.Item(38, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.295'
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.40%  '
$cell = $ws.Cells.Item(39, 2)
$cell.NumberFormat = "@"
$cell.Value = 'VeChain'
$cell = $ws.Cells.Item(39, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.01709'
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -4.60%  '
$cell = $ws.Cells.Item(40, 2)
$cell.NumberFormat = "@"
$cell.Value = 'TheSandbox'
$cell = $ws.Cells.Item(40, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4718'
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.05%  '
$cell = $ws.Cells.Item(41, 2)
$cell.NumberFormat = "@"
$cell.Value = 'TrustWalletToken'
$cell = $ws.Cells.Item(41, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9053'
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.74%  '
$cell = $ws.Cells.Item(42, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Quant'
$cell = $ws.Cells.Item(42, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '107.49'
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.21%  '
$cell = $ws.Cells.Item(43, 2)
$cell.NumberFormat = "@"
$cell.Value = 'FraxShare'
$cell = $ws.Cells.Item(43, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.876'
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.47%  '
$cell = $ws.Cells.Item(44, 2)
$cell.NumberFormat = "@"
$cell.Value = 'PaxDollar'
$cell = $ws.Cells.Item(44, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.17%  '
$cell = $ws.Cells.Item(45, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Aptos'
$cell = $ws.Cells.Item(45, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.360'
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.67%  '
$cell = $ws.Cells.Item(46, 2)
$cell.NumberFormat = "@"
$cell.Value = 'EnergySwap'
$cell = $ws.Cells.Item(46, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.967'
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.21%  '
$cell = $ws.Cells.Item(47, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Decentraland'
$cell = $ws.Cells.Item(47, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4078'
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.16%  '
$cell = $ws.Cells.Item(48, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Algorand'
$cell = $ws.Cells.Item(48, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1231'
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.79%  '
$cell = $ws.Cells.Item(49, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Elrond'
$cell = $ws.Cells.Item(49, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '34.79'
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.64%  '
$cell = $ws.Cells.Item(50, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Cronos'
$cell = $ws.Cells.Item(50, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05792'
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.61%  '
$cell = $ws.Cells.Item(51, 2)
$cell.NumberFormat = "@"
$cell.Value = 'EOS'
$cell = $ws.Cells.Item(51, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.8915'
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.20%  '
